$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.989.03'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '1.559.33'
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +0.98%  '
$ws.Range("E10").Value = '  +1.84%  '
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").Value = '1.563.45'
$ws.Range("E13").Value = '  +0.76%  '
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("D16").Value = '27.007.04'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '0.0₃0705'
$ws.Range("E18").Value = '  +1.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("E20").Value = '  +1.45%  '
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.25%  '
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("E24").Value = '  -0.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.00%  '
$ws.Range("E26").Value = '  -0.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.74%  '
$ws.Range("E28").Value = '  +1.59%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  +1.44%  '
$ws.Range("E31").Value = '  +3.37%  '
$ws.Range("E32").Value = '  +0.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.19'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.02%  '
$ws.Range("D34").Value = '1.422.99'
$ws.Range("E34").Value = '  +0.30%  '
$ws.Range("E35").Value = '  +11.28%  '
$ws.Range("E36").Value = '  +1.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.30%  '
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("E39").Value = '  +2.31%  '
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("E41").Value = '  +0.46%  '
$ws.Range("E42").Value = '  +0.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.61%  '
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("E46").Value = '  -0.93%  '
$ws.Range("D47").Value = '1.696.24'
$ws.Range("E47").Value = '  +0.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.01%  '
$ws.Range("E49").Value = '  +2.87%  '
$ws.Range("E50").Value = '  -0.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0958'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.81%  '
